$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 293, pushing the existing rows 293:404 down to 294:405
$ws.Rows("293:293").Insert()

# Populate the newly inserted row 293 with the new record.
# Columns A,B,C,E,F,G,H,I,J,K,L,Q,T repeat the constant values used throughout
# this data set; D,M,N,O,P,R,S hold the new record's unique values.
$ws.Cells.Item(293, 1).Value2 = 10
$ws.Cells.Item(293, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(293, 3).Value2 = "La Araucanía"
$ws.Cells.Item(293, 4).Value2 = 44795
$ws.Cells.Item(293, 5).Value2 = 9
$ws.Cells.Item(293, 6).Value2 = "Fruta"
$ws.Cells.Item(293, 7).Value2 = 100108
$ws.Cells.Item(293, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(293, 9).Value2 = 100108002
$ws.Cells.Item(293, 10).Value2 = "Mango"
$ws.Cells.Item(293, 11).Value2 = "Sin especificar"
$ws.Cells.Item(293, 12).Value2 = "Primera"
$ws.Cells.Item(293, 13).Value2 = 125
$ws.Cells.Item(293, 14).Value2 = 12000
$ws.Cells.Item(293, 15).Value2 = 12000
$ws.Cells.Item(293, 16).Value2 = 12000
$ws.Cells.Item(293, 17).Value2 = "$/bandeja 4 kilos"
$ws.Cells.Item(293, 18).Value2 = "Brasil"
$ws.Cells.Item(293, 19).Value2 = 3000
$ws.Cells.Item(293, 20).Value2 = 4

# Make sure the date cell keeps the existing date number format (style carries
# over from the row-insert, but set it explicitly to be safe).
$ws.Cells.Item(293, 4).NumberFormat = $ws.Cells.Item(294, 4).NumberFormat
